# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Halicarnassus Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 448
$ws.Range("I98").Value = 410.8
$ws.Range("K98").Value = 410.8
$ws.Range("M98").Value = 1087.2
$ws.Range("H99").Value = 2487.2
$ws.Range("I99").Value = 158.66667
$ws.Range("J99").Value = 5980
$ws.Range("K99").Value = 476.00001
$ws.Range("L99").Value = 17940
$ws.Range("M99").Value = 1021.99999
$ws.Range("N99").Value = -20936
$ws.Range("H113").Value = 4317.5
$ws.Range("I113").Value = 3681
$ws.Range("K113").Value = 3681
$ws.Range("M113").Value = -427
$ws.Range("H122").Value = 448
$ws.Range("I122").Value = 410.8
$ws.Range("K122").Value = 1232.4
$ws.Range("M122").Value = 1217.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 2429.1667
$ws.Range("I36").Value = 1694
$ws.Range("J36").Value = 3899.5
$ws.Range("K36").Value = 1694
$ws.Range("L36").Value = 3899.5
$ws.Range("M36").Value = -1348
$ws.Range("N36").Value = -4591.5
$ws.Range("H61").Value = 1802.2632
$ws.Range("I61").Value = 1259.3125
$ws.Range("K61").Value = 1259.3125
$ws.Range("M61").Value = -1047.3125
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H124").Value = 25723.25
$ws.Range("J124").Value = 25723.25
$ws.Range("L124").Value = 25723.25
$ws.Range("N124").Value = -35543.25
$ws.Range("H132").Value = 1587.5
$ws.Range("I132").Value = 1296.875
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 3890.625
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -1360.625
$ws.Range("N132").Value = -13310
$ws.Range("H136").Value = 1802.2632
$ws.Range("I136").Value = 1259.3125
$ws.Range("K136").Value = 3777.9375
$ws.Range("M136").Value = -1227.9375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H105").Value = 6163389.5
$ws.Range("I105").Value = 10084228
$ws.Range("K105").Value = 10084228
$ws.Range("M105").Value = -10082481
$ws.Range("H132").Value = 99996.664
$ws.Range("J132").Value = 99996.664
$ws.Range("L132").Value = 99996.664
$ws.Range("N132").Value = -110116.664
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1049.5454
$ws.Range("I32").Value = 1054.5
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 1054.5
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -738.5
$ws.Range("N32").Value = -1632
$ws.Range("H99").Value = 2800.3333
$ws.Range("I99").Value = 2975.75
$ws.Range("K99").Value = 2975.75
$ws.Range("M99").Value = -1477.75
$ws.Range("H126").Value = 2800.3333
$ws.Range("I126").Value = 2975.75
$ws.Range("K126").Value = 8927.25
$ws.Range("M126").Value = -6457.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1972.5333
$ws.Range("I5").Value = 2018
$ws.Range("J5").Value = 1949.8
$ws.Range("K5").Value = 6054
$ws.Range("L5").Value = 5849.4
$ws.Range("M5").Value = -5942
$ws.Range("N5").Value = -6073.4
$ws.Range("H132").Value = 2890.8572
$ws.Range("I132").Value = 1873.1666
$ws.Range("J132").Value = 3654.125
$ws.Range("K132").Value = 16858.4994
$ws.Range("L132").Value = 32887.125
$ws.Range("M132").Value = -14328.4994
$ws.Range("N132").Value = -37947.125
$ws.Range("H134").Value = 2133.2222
$ws.Range("I134").Value = 1885.5714
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5656.7142
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -586.7142000000003
$ws.Range("N134").Value = -19140
$ws.Range("H135").Value = 1972.5333
$ws.Range("I135").Value = 2018
$ws.Range("J135").Value = 1949.8
$ws.Range("K135").Value = 18162
$ws.Range("L135").Value = 17548.2
$ws.Range("M135").Value = -15627
$ws.Range("N135").Value = -22618.2
$ws.Range("H136").Value = 13000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H139").Value = 2740
$ws.Range("I139").Value = 480
$ws.Range("K139").Value = 1440
$ws.Range("M139").Value = 3700
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 69696
$ws.Range("J57").Value = 69696
$ws.Range("L57").Value = 69696
$ws.Range("N57").Value = -71336
$ws.Range("H100").Value = 39999
$ws.Range("J100").Value = 39999
$ws.Range("L100").Value = 39999
$ws.Range("N100").Value = -42163
$ws.Range("H122").Value = 1681.4
$ws.Range("I122").Value = 956.5454999999999
$ws.Range("K122").Value = 2869.6365
$ws.Range("M122").Value = -419.6364999999996
$ws.Range("H132").Value = 3882.5
$ws.Range("I132").Value = 4294.737
$ws.Range("J132").Value = 1271.6666
$ws.Range("K132").Value = 12884.211
$ws.Range("L132").Value = 3814.9998
$ws.Range("M132").Value = -10354.211
$ws.Range("N132").Value = -8874.9998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 6867.2
$ws.Range("J4").Value = 567
$ws.Range("L4").Value = 567
$ws.Range("N4").Value = -793
$ws.Range("H132").Value = 2956.3635
$ws.Range("I132").Value = 2946.6667
$ws.Range("K132").Value = 8840.000100000001
$ws.Range("M132").Value = -6310.000100000001
$ws.Range("H136").Value = 2571.88
$ws.Range("I136").Value = 1748
$ws.Range("K136").Value = 5244
$ws.Range("M136").Value = -2694
